$p = $ppt.ActivePresentation

# 1) Move the vertical slide guide from 5280 -> 5496 (1/100 mm units... actually points)
$guides = $p.Guides
for ($i = 1; $i -le $guides.Count; $i++) {
    $g = $guides.Item($i)
    if ($g.Position -eq 5280) {
        $g.Position = 5496
    }
}
